# New weekly price-report row arrived: insert it at the top of the
# "Vega Modelo de Temuco" Berenjena block (row 245) and push the
# existing rows (245-306) down by one, creating a new last row 307.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(245).Insert()

$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 44754
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = 100112001
$ws.Range("G245").Value = "Berenjena"
$ws.Range("H245").Value = "Sin especificar"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 100
$ws.Range("K245").Value = 13000
$ws.Range("L245").Value = 14000
$ws.Range("M245").Value = 13400
$ws.Range("N245").Value = "`$/caja 60 unidades"
$ws.Range("O245").Value = "Región de Arica y Parinacota"
$ws.Range("P245").Value = 223
$ws.Range("Q245").Value = 60
$ws.Range("R245").Value = "Hortaliza"
